$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Paragraph "bazwilks changes" -> split into runs with proofErr tags
#    around "bazwilks": <proofErr spellStart/><r>bazwilks</r><proofErr
#    spellEnd/><r> changes</r>
# ----------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "bazwilks changes") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Range covering only the paragraph's content, excluding the paragraph mark,
    # so InsertXML replaces the runs but keeps the paragraph's own identity/mark.
    $innerEnd = $target.Range.End - 1
    $inner = $d.Range($target.Range.Start, $innerEnd)

    $xmlBazwilks = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>bazwilks</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> changes</w:t></w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $inner.InsertXML($xmlBazwilks) | Out-Null
}

# ----------------------------------------------------------------------
# 2) After the paragraph "Test - Thursday 30 July" add three new
#    paragraphs:
#      a) "Test - Thursday " + "15 March" (two runs)
#      b) "Test - Thursday 15 March" (one run)
#      c) "Test - Thursday 15 March" (one run)
# ----------------------------------------------------------------------
$afterPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Test – Thursday 30 July") {
        $afterPara = $p
        break
    }
}

if ($afterPara -ne $null) {
    $insertAt = $d.Range($afterPara.Range.End, $afterPara.Range.End)

    # A trailing empty <w:p/> is appended deliberately: this engine's
    # InsertXML merges the *last* <w:p> in the fragment with whatever
    # paragraph sits at the insertion point (absorbing its paragraph
    # mark/identity). Supplying a throwaway empty paragraph at the end
    # means that merge "eats" the dummy instead of the real following
    # paragraph, so the three real paragraphs end up as genuinely new,
    # distinct paragraphs and the original following paragraph keeps
    # its identity untouched. The dummy paragraph is removed afterwards.
    $xmlNewParas = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:r><w:t xml:space="preserve">Test – Thursday </w:t></w:r><w:r><w:t>15 March</w:t></w:r></w:p>' +
        '<w:p><w:r><w:t>Test – Thursday 15 March</w:t></w:r></w:p>' +
        '<w:p><w:r><w:t>Test – Thursday 15 March</w:t></w:r></w:p>' +
        '<w:p/>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $countBefore = $d.Paragraphs.Count
    $insertAt.InsertXML($xmlNewParas) | Out-Null
    $countAfter = $d.Paragraphs.Count
    $added = $countAfter - $countBefore

    # The dummy paragraph that absorbed the merge is the one right
    # after our three newly-inserted paragraphs (i.e. the last of the
    # freshly added ones). Find $afterPara's current index by position,
    # then step forward by the number of paragraphs that were added.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Start -eq $afterPara.Range.Start) {
            $dummyIndex = $i + $added
            break
        }
    }
    $dummy = $d.Paragraphs($dummyIndex)
    $dummy.Range.Delete()
}
